$wb = $excel.ActiveWorkbook

# ALC!row87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 52631
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 52631
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 52631
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -55127

# ALC!row90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 52631
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 52631
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 157893
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -170373

# ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 7248016
$ws.Range("I100").Value = 10418254
$ws.Range("J100").Value = 1756.5714
$ws.Range("K100").Value = 10418254
$ws.Range("L100").Value = 1756.5714
$ws.Range("M100").Value = -10417713
$ws.Range("N100").Value = -2838.5714

# ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 13637652
$ws.Range("J112").Value = 14355371
$ws.Range("L112").Value = 43066113
$ws.Range("N112").Value = -43068329

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1578
$ws.Range("I129").Value = 493.33334
$ws.Range("J129").Value = 2042.8572
$ws.Range("K129").Value = 1480.00002
$ws.Range("L129").Value = 6128.571599999999
$ws.Range("M129").Value = 3519.99998
$ws.Range("N129").Value = -16128.5716

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2069.65
$ws.Range("I138").Value = 1002.8205
$ws.Range("J138").Value = 2751.7212
$ws.Range("K138").Value = 3008.4615
$ws.Range("L138").Value = 8255.1636
$ws.Range("M138").Value = 2131.5385
$ws.Range("N138").Value = -18535.1636

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6295.75
$ws.Range("I32").Value = 4288.811
$ws.Range("J32").Value = 24358.2
$ws.Range("K32").Value = 4288.811
$ws.Range("L32").Value = 24358.2
$ws.Range("M32").Value = -4001.811
$ws.Range("N32").Value = -24932.2

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3000.4546
$ws.Range("I61").Value = 2044.5333
$ws.Range("J61").Value = 5048.857
$ws.Range("K61").Value = 2044.5333
$ws.Range("L61").Value = 5048.857
$ws.Range("M61").Value = -1832.5333
$ws.Range("N61").Value = -5472.857

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12035.5
$ws.Range("I74").Value = 2062.25
$ws.Range("J74").Value = 31982
$ws.Range("K74").Value = 2062.25
$ws.Range("L74").Value = 31982
$ws.Range("M74").Value = -1188.25
$ws.Range("N74").Value = -33730

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 12035.5
$ws.Range("I77").Value = 2062.25
$ws.Range("J77").Value = 31982
$ws.Range("K77").Value = 10311.25
$ws.Range("L77").Value = 159910
$ws.Range("M77").Value = -5943.25
$ws.Range("N77").Value = -168646

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 7836.357
$ws.Range("I97").Value = 10571
$ws.Range("J97").Value = 999.75
$ws.Range("K97").Value = 10571
$ws.Range("L97").Value = 999.75
$ws.Range("M97").Value = -10075
$ws.Range("N97").Value = -1991.75

# ARM!row104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3000.4546
$ws.Range("I136").Value = 2044.5333
$ws.Range("J136").Value = 5048.857
$ws.Range("K136").Value = 6133.5999
$ws.Range("L136").Value = 15146.571
$ws.Range("M136").Value = -3583.5999
$ws.Range("N136").Value = -20246.571

# BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 202.0625
$ws.Range("I22").Value = 205.61539
$ws.Range("J22").Value = 186.66667
$ws.Range("K22").Value = 205.61539
$ws.Range("L22").Value = 186.66667
$ws.Range("M22").Value = -32.61538999999999
$ws.Range("N22").Value = -532.6666700000001

# BSM!row58
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 21300
$ws.Range("J58").Value = 21300
$ws.Range("L58").Value = 21300
$ws.Range("N58").Value = -21888

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3359.2856
$ws.Range("I94").Value = 2657.2727
$ws.Range("J94").Value = 5933.3335
$ws.Range("K94").Value = 2657.2727
$ws.Range("L94").Value = 5933.3335
$ws.Range("M94").Value = -2206.2727
$ws.Range("N94").Value = -6835.3335

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3745.1875
$ws.Range("I134").Value = 2526.5264
$ws.Range("J134").Value = 5526.3076
$ws.Range("K134").Value = 7579.5792
$ws.Range("L134").Value = 16578.9228
$ws.Range("M134").Value = -5044.5792
$ws.Range("N134").Value = -21648.9228

# CRP!row7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3730.4092
$ws.Range("I31").Value = 1649.5151
$ws.Range("J31").Value = 5811.303
$ws.Range("K31").Value = 1649.5151
$ws.Range("L31").Value = 5811.303
$ws.Range("M31").Value = -1354.5151
$ws.Range("N31").Value = -6401.303

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3730.4092
$ws.Range("I34").Value = 1649.5151
$ws.Range("J34").Value = 5811.303
$ws.Range("K34").Value = 1649.5151
$ws.Range("L34").Value = 5811.303
$ws.Range("M34").Value = -1447.5151
$ws.Range("N34").Value = -6215.303

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 32259506
$ws.Range("I58").Value = 35715364
$ws.Range("J58").Value = 4811.3335
$ws.Range("K58").Value = 35715364
$ws.Range("L58").Value = 4811.3335
$ws.Range("M58").Value = -35715161
$ws.Range("N58").Value = -5217.3335

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 34093504
$ws.Range("I134").Value = 40001580
$ws.Range("J134").Value = 26319720
$ws.Range("K134").Value = 120004740
$ws.Range("L134").Value = 78959160
$ws.Range("M134").Value = -120002205
$ws.Range("N134").Value = -78964230

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 32259506
$ws.Range("I136").Value = 35715364
$ws.Range("J136").Value = 4811.3335
$ws.Range("K136").Value = 107146092
$ws.Range("L136").Value = 14434.0005
$ws.Range("M136").Value = -107143542
$ws.Range("N136").Value = -19534.0005

# CUL!row33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 64.14286
$ws.Range("I33").Value = 62.923077
$ws.Range("J33").Value = 80
$ws.Range("K33").Value = 377.538462
$ws.Range("L33").Value = 480
$ws.Range("M33").Value = -94.53846199999998
$ws.Range("N33").Value = -1046

# CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4646.3335
$ws.Range("I34").Value = 400
$ws.Range("J34").Value = 4949.643
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 14848.929
$ws.Range("M34").Value = -1116
$ws.Range("N34").Value = -15016.929

# CUL!row46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 929.2308
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 1058
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 3174
$ws.Range("M46").Value = -1409
$ws.Range("N46").Value = -3356

# CUL!row109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 920.2
$ws.Range("I109").Value = 920.2
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 2760.6
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -1720.6
$ws.Range("N109").ClearContents()

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3637101.2
$ws.Range("I113").Value = 636.6
$ws.Range("J113").Value = 6061411
$ws.Range("K113").Value = 1909.8
$ws.Range("L113").Value = 18184233
$ws.Range("M113").Value = 260.1999999999998
$ws.Range("N113").Value = -18188573

# CUL!row121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 903.4706
$ws.Range("I121").Value = 251.11111
$ws.Range("J121").Value = 1637.375
$ws.Range("K121").Value = 753.3333299999999
$ws.Range("L121").Value = 4912.125
$ws.Range("M121").Value = 556.6666700000001
$ws.Range("N121").Value = -7532.125

# CUL!row139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 7089.5654
$ws.Range("I139").Value = 5537.3335
$ws.Range("J139").Value = 10000
$ws.Range("K139").Value = 16612.0005
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -11472.0005
$ws.Range("N139").Value = -40280

# GSM!row7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1500250
$ws.Range("I7").Value = 1500250
$ws.Range("K7").Value = 1500250
$ws.Range("M7").Value = -1500138

# GSM!row8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 1500250
$ws.Range("I8").Value = 1500250
$ws.Range("K8").Value = 1500250
$ws.Range("M8").Value = -1500111

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 925.1429000000001
$ws.Range("I122").Value = 925.1429000000001
$ws.Range("K122").Value = 2775.4287
$ws.Range("M122").Value = -325.4287000000004

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2946.4707
$ws.Range("I126").Value = 2834.5454
$ws.Range("K126").Value = 8503.636200000001
$ws.Range("M126").Value = -6033.636200000001

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2809.5667
$ws.Range("J7").Value = 3342.3684
$ws.Range("L7").Value = 3342.3684
$ws.Range("N7").Value = -3566.3684

# LTW!row93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1832.125
$ws.Range("I93").Value = 1583.3572
$ws.Range("J93").Value = 2180.4
$ws.Range("K93").Value = 1583.3572
$ws.Range("L93").Value = 2180.4
$ws.Range("M93").Value = -335.3571999999999
$ws.Range("N93").Value = -4676.4

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3860
$ws.Range("I122").Value = 3580
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10740
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8290
$ws.Range("N122").Value = -16900

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2809.5667
$ws.Range("J126").Value = 3342.3684
$ws.Range("L126").Value = 10027.1052
$ws.Range("N126").Value = -14967.1052

# LTW!row127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 47866.668
$ws.Range("J127").Value = 47866.668
$ws.Range("L127").Value = 47866.668
$ws.Range("N127").Value = -57786.668

# WVR!row14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 17000
$ws.Range("I14").Value = 26250
$ws.Range("J14").Value = 7750
$ws.Range("K14").Value = 26250
$ws.Range("L14").Value = 7750
$ws.Range("M14").Value = -26082
$ws.Range("N14").Value = -8086

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 837621
$ws.Range("I81").Value = 3335517.2
$ws.Range("J81").Value = 4988.8887
$ws.Range("K81").Value = 6671034.4
$ws.Range("L81").Value = 9977.777400000001
$ws.Range("M81").Value = -6669973.4
$ws.Range("N81").Value = -12099.7774

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 837621
$ws.Range("I84").Value = 3335517.2
$ws.Range("J84").Value = 4988.8887
$ws.Range("K84").Value = 33355172
$ws.Range("L84").Value = 49888.887
$ws.Range("M84").Value = -33349868
$ws.Range("N84").Value = -60496.887

# WVR!row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 66668284
$ws.Range("I96").Value = 90910440
$ws.Range("J96").Value = 2361
$ws.Range("K96").Value = 90910440
$ws.Range("L96").Value = 2361
$ws.Range("M96").Value = -90909067
$ws.Range("N96").Value = -5107

# WVR!row109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 33500
$ws.Range("J109").Value = 33500
$ws.Range("L109").Value = 33500
$ws.Range("N109").Value = -36274

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3291.5652
$ws.Range("I132").Value = 3158.75
$ws.Range("K132").Value = 9476.25
$ws.Range("M132").Value = -6946.25
$ws.Range("N132").Value = -15946.25
